$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.676
$ws.Range("D6").Value = -8.169
$ws.Range("D7").Value = -7.585000000000001
$ws.Range("E7").Value = 16.852
$ws.Range("E12").Value = 17.638
$ws.Range("E15").Value = 15.955
$ws.Range("D16").Value = -8.626000000000001
$ws.Range("D20").Value = -7.808000000000002
$ws.Range("E20").Value = 16.457
$ws.Range("E21").Value = 16.523
$ws.Range("E22").Value = 16.506
$ws.Range("E23").Value = 16.434
$ws.Range("D28").Value = -8.055
$ws.Range("D29").Value = -7.56
$ws.Range("E29").Value = 16.689
$ws.Range("D32").Value = -7.98
$ws.Range("E34").Value = 16.846
$ws.Range("D40").Value = -7.883999999999999
$ws.Range("E42").Value = 16.539
$ws.Range("E43").Value = 16.951
$ws.Range("E44").Value = 16.518
$ws.Range("E45").Value = 16.774
$ws.Range("D46").Value = -8.022
$ws.Range("E46").Value = 16.968
$ws.Range("E50").Value = 16.438
$ws.Range("D51").Value = -8.241000000000001
$ws.Range("E51").Value = 16.489
$ws.Range("D52").Value = -7.87
$ws.Range("D57").Value = -8.06
$ws.Range("D59").Value = -8.123999999999999
$ws.Range("D62").Value = -7.946000000000001
$ws.Range("D66").Value = -7.419
$ws.Range("E66").Value = 17.325
$ws.Range("E67").Value = 17.241
$ws.Range("D73").Value = -7.789
$ws.Range("D74").Value = -7.904999999999998
$ws.Range("E79").Value = 16.957
$ws.Range("E84").Value = 16.606
$ws.Range("D92").Value = -7.676
$ws.Range("E92").Value = 16.832
$ws.Range("E97").Value = 16.787
$ws.Range("D100").Value = -8.295
